# Update the "Price" (D) and "Volume(1h)" (E) columns with freshly scraped
# values, plus two coin-rank swaps (rows 11/12 and 36/37). Every cell in this
# sheet is stored as plain text (the source file is produced by openpyxl, not
# real Excel), so Price values that look numeric are written with a leading
# apostrophe to force Excel to keep them as text instead of auto-converting
# them to numbers; the apostrophe itself is not stored in the cell value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '40.961.93'
$ws.Range('E2').Value = '  -2.10%  '

# Row 3
$ws.Range('D3').Value = '2.171.24'
$ws.Range('E3').Value = '  -2.98%  '

# Row 4
$ws.Range('E4').Value = '  -0.08%  '

# Row 5
$ws.Range('D5').Value = '''248.45'
$ws.Range('E5').Value = '  -1.25%  '

# Row 6
$ws.Range('D6').Value = '''0.618'
$ws.Range('E6').Value = '  -1.63%  '

# Row 7
$ws.Range('D7').Value = '''66.91'
$ws.Range('E7').Value = '  -7.22%  '

# Row 8
$ws.Range('E8').Value = '  +0.03%  '

# Row 9
$ws.Range('D9').Value = '''0.564'
$ws.Range('E9').Value = '  -0.24%  '

# Row 10
$ws.Range('D10').Value = '''58.31'
$ws.Range('E10').Value = '  -0.16%  '

# Row 11
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').Value = '''35.75'
$ws.Range('E11').Value = '  -15.30%  '

# Row 12
$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').Value = '''0.0924'
$ws.Range('E12').Value = '  -5.04%  '

# Row 13
$ws.Range('E13').Value = '  -1.61%  '

# Row 14
$ws.Range('D14').Value = '''6.89'
$ws.Range('E14').Value = '  -0.17%  '

# Row 15
$ws.Range('D15').Value = '2.495.67'
$ws.Range('E15').Value = '  -3.10%  '

# Row 16
$ws.Range('D16').Value = '''0.860'

# Row 17
$ws.Range('E17').Value = '  -6.66%  '

# Row 18
$ws.Range('D18').Value = '2.183.68'
$ws.Range('E18').Value = '  -2.64%  '

# Row 19
$ws.Range('D19').Value = '40.824.45'
$ws.Range('E19').Value = '  -2.45%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0937'
$ws.Range('E20').Value = '  -3.19%  '

# Row 21
$ws.Range('D21').Value = '''6.07'
$ws.Range('E21').Value = '  -2.72%  '

# Row 22
$ws.Range('D22').Value = '''71.32'
$ws.Range('E22').Value = '  -2.86%  '

# Row 23
$ws.Range('D23').Value = '''229.49'
$ws.Range('E23').Value = '  -2.60%  '

# Row 24
$ws.Range('D24').Value = '''2.06'
$ws.Range('E24').Value = '  -8.41%  '

# Row 25
$ws.Range('D25').Value = '''1.00'
$ws.Range('E25').Value = '  +0.02%  '

# Row 26
$ws.Range('E26').Value = '  +11.77%  '

# Row 27
$ws.Range('D27').Value = '''3.73'
$ws.Range('E27').Value = '  -0.77%  '

# Row 28
$ws.Range('D28').Value = '''2.41'
$ws.Range('E28').Value = '  -3.80%  '

# Row 29
$ws.Range('D29').Value = '''2.13'
$ws.Range('E29').Value = '  -3.20%  '

# Row 30
$ws.Range('D30').Value = '''167.77'
$ws.Range('E30').Value = '  -2.56%  '

# Row 31
$ws.Range('D31').Value = '''20.16'
$ws.Range('E31').Value = '  -3.07%  '

# Row 32
$ws.Range('D32').Value = '''0.121'
$ws.Range('E32').Value = '  -1.52%  '

# Row 33
$ws.Range('D33').Value = '''5.62'
$ws.Range('E33').Value = '  +2.63%  '

# Row 34
$ws.Range('D34').Value = '''0.0742'
$ws.Range('E34').Value = '  +2.66%  '

# Row 35
$ws.Range('D35').Value = '''0.121'
$ws.Range('E35').Value = '  -3.34%  '

# Row 36
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').Value = '''4.54'
$ws.Range('E36').Value = '  -3.66%  '

# Row 37
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = '''4.08'
$ws.Range('E37').Value = '  -2.04%  '

# Row 38
$ws.Range('D38').Value = '''25.20'
$ws.Range('E38').Value = '  -4.67%  '

# Row 39
$ws.Range('E39').Value = '  +7.20%  '

# Row 40
$ws.Range('D40').Value = '''2.17'
$ws.Range('E40').Value = '  -5.37%  '

# Row 41
$ws.Range('D41').Value = '''5.49'
$ws.Range('E41').Value = '  -9.40%  '

# Row 42
$ws.Range('D42').Value = '''11.58'
$ws.Range('E42').Value = '  -0.96%  '

# Row 43
$ws.Range('D43').Value = '''61.21'
$ws.Range('E43').Value = '  -10.84%  '

# Row 44
$ws.Range('D44').Value = '''4.73'
$ws.Range('E44').Value = '  -6.21%  '

# Row 45
$ws.Range('D45').Value = '''0.192'
$ws.Range('E45').Value = '  -11.15%  '

# Row 46
$ws.Range('D46').Value = '''8.52'
$ws.Range('E46').Value = '  -3.73%  '

# Row 47
$ws.Range('D47').Value = '''1.01'
$ws.Range('E47').Value = '  +0.54%  '

# Row 48
$ws.Range('E48').Value = '  +4.44%  '

# Row 49
$ws.Range('D49').Value = '''0.0981'
$ws.Range('E49').Value = '  -3.72%  '

# Row 50
$ws.Range('D50').Value = '''1.14'
$ws.Range('E50').Value = '  -4.14%  '

# Row 51
$ws.Range('D51').Value = '''2.70'
$ws.Range('E51').Value = '  -0.70%  '
